$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 53, shifting existing rows 53-67 down to 54-68
$ws.Rows.Item(53).Insert()

# Populate the newly inserted row 53 with the new record's data
$ws.Range("A53").Value = 5
$ws.Range("B53").Value = 'Macroferia Regional de Talca'
$ws.Range("C53").Value = 'Maule'
$ws.Range("D53").Value = 44468
$ws.Range("E53").Value = 7
$ws.Range("F53").Value = 100112001
$ws.Range("G53").Value = 'Berenjena'
$ws.Range("H53").Value = 'Sin especificar'
$ws.Range("I53").Value = 'Primera'
$ws.Range("J53").Value = 200
$ws.Range("K53").Value = 8000
$ws.Range("L53").Value = 8000
$ws.Range("M53").Value = 8000
$ws.Range("N53").Value = '$/caja 50 unidades'
$ws.Range("O53").Value = 'Región de Arica y Parinacota'
$ws.Range("P53").Value = 160
$ws.Range("Q53").Value = 50
$ws.Range("R53").Value = 'Hortaliza'
